# Them option them giao vien
# Update existing two student rows and add a new row for a new student/teacher entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: update "So hieu" (kept as text because of its number-format style) and the name
$ws.Range("A2").Value = "49.323"
$ws.Range("D2").Value = "Bui Quang Tuan "
$ws.Range("E2").Value = "B12D49"

# Row 3: update numeric id, name, and class
$ws.Range("A3").Value = 49.323999999999998
$ws.Range("D3").Value = "Nguyen Anh Tuan "
$ws.Range("E3").Value = "B12D49"
$ws.Range("F3").Value = 3

# Row 4 (new): add a new entry
$ws.Range("A4").Value = 49.325000000000003
$ws.Range("B4").Value = 2019
$ws.Range("C4").Value = 2020
$ws.Range("D4").Value = "Nguyen Thi Lam Vien"
$ws.Range("E4").Value = "B13D49"
$ws.Range("F4").Value = 4

# Move the active selection to match the new last empty-ish cell
$ws.Range("G4").Select()
